$d = $word.ActiveDocument

# The document ends with a list paragraph "...Добавить testpoints на все
# сигналы." that also carries the (collapsed) "_GoBack" bookmark right
# before its paragraph mark. We need to split that paragraph in two,
# appending a brand-new list item "Переразвести цепи в соответствии со
# схемой." and keep the "_GoBack" bookmark trailing the new last
# paragraph (same relative place it held before).
#
# Inserting a paragraph mark straight through the bookmark's own Range
# relocates the bookmark to the wrong place in this runtime, so instead
# we: 1) append the new sentence using the bookmark's Range (so the
# bookmark naturally keeps trailing the appended text), and 2) locate the
# boundary between the old and the new sentence via Find and insert the
# paragraph break there.

$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertBefore("Переразвести цепи в соответствии со схемой.")

$rng = $d.Content
$rng.Find.Execute(" на все сигналы.", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null
